$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2675.75
$ws.Range("J17").Value = 2675.75
$ws.Range("L17").Value = 8027.25
$ws.Range("N17").Value = -8363.25
$ws.Range("H86").Value = 6581733
$ws.Range("I86").Value = 4662.3335
$ws.Range("J86").Value = 10527976
$ws.Range("K86").Value = 4662.3335
$ws.Range("L86").Value = 10527976
$ws.Range("M86").Value = -3539.3335
$ws.Range("N86").Value = -10530222
$ws.Range("H89").Value = 6581733
$ws.Range("I89").Value = 4662.3335
$ws.Range("J89").Value = 10527976
$ws.Range("K89").Value = 23311.6675
$ws.Range("L89").Value = 52639880
$ws.Range("M89").Value = -17695.6675
$ws.Range("N89").Value = -52651112
$ws.Range("H98").Value = 1582.0851
$ws.Range("I98").Value = 1159.7894
$ws.Range("J98").Value = 3365.111
$ws.Range("K98").Value = 1159.7894
$ws.Range("L98").Value = 3365.111
$ws.Range("M98").Value = 338.2106000000001
$ws.Range("N98").Value = -6361.111
$ws.Range("H100").Value = 10880.929
$ws.Range("I100").Value = 677.5
$ws.Range("J100").Value = 12581.5
$ws.Range("K100").Value = 677.5
$ws.Range("L100").Value = 12581.5
$ws.Range("M100").Value = -136.5
$ws.Range("N100").Value = -13663.5
$ws.Range("H122").Value = 1582.0851
$ws.Range("I122").Value = 1159.7894
$ws.Range("J122").Value = 3365.111
$ws.Range("K122").Value = 3479.3682
$ws.Range("L122").Value = 10095.333
$ws.Range("M122").Value = -1029.3682
$ws.Range("N122").Value = -14995.333
$ws.Range("H123").Value = 75226.63
$ws.Range("J123").Value = 75226.63
$ws.Range("L123").Value = 75226.63
$ws.Range("N123").Value = -85026.63
$ws.Range("H138").Value = 4960.7817
$ws.Range("I138").Value = 2863.8462
$ws.Range("J138").Value = 5854.5576
$ws.Range("K138").Value = 8591.5386
$ws.Range("L138").Value = 17563.6728
$ws.Range("M138").Value = -3451.5386
$ws.Range("N138").Value = -27843.6728
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 2263.3555
$ws.Range("I32").Value = 1442.2289
$ws.Range("K32").Value = 1442.2289
$ws.Range("M32").Value = -1155.2289
$ws.Range("H122").Value = 3215.6943
$ws.Range("I122").Value = 1834.5625
$ws.Range("J122").Value = 4320.6
$ws.Range("K122").Value = 5503.6875
$ws.Range("L122").Value = 12961.8
$ws.Range("M122").Value = -3053.6875
$ws.Range("N122").Value = -17861.8
$ws.Range("H132").Value = 2463.46
$ws.Range("J132").Value = 9428.571
$ws.Range("L132").Value = 28285.713
$ws.Range("N132").Value = -33345.713
$ws.Range("H141").Value = 54999.5
$ws.Range("J141").Value = 54999.5
$ws.Range("L141").Value = 54999.5
$ws.Range("N141").Value = -65359.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1132.75
$ws.Range("I20").Value = 744.0625
$ws.Range("J20").Value = 2687.5
$ws.Range("K20").Value = 744.0625
$ws.Range("L20").Value = 2687.5
$ws.Range("M20").Value = -497.0625
$ws.Range("N20").Value = -3181.5
$ws.Range("H50").Value = 48576
$ws.Range("J50").Value = 48576
$ws.Range("L50").Value = 48576
$ws.Range("N50").Value = -49724
$ws.Range("H64").Value = 1549
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 1565.3334
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 1565.3334
$ws.Range("M64").Value = -1275
$ws.Range("N64").Value = -2015.3334
$ws.Range("H67").Value = 1549
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 1565.3334
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 1565.3334
$ws.Range("M67").Value = -720
$ws.Range("N67").Value = -3125.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 361489.4
$ws.Range("J31").Value = 7692.846
$ws.Range("L31").Value = 7692.846
$ws.Range("N31").Value = -8282.846
$ws.Range("H34").Value = 361489.4
$ws.Range("J34").Value = 7692.846
$ws.Range("L34").Value = 7692.846
$ws.Range("N34").Value = -8096.846
$ws.Range("H58").Value = 225959.89
$ws.Range("I58").Value = 557207.25
$ws.Range("K58").Value = 557207.25
$ws.Range("M58").Value = -557004.25
$ws.Range("H136").Value = 225959.89
$ws.Range("I136").Value = 557207.25
$ws.Range("K136").Value = 1671621.75
$ws.Range("M136").Value = -1669071.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2888
$ws.Range("H68").Value = 2859938
$ws.Range("J68").Value = 2002580
$ws.Range("L68").Value = 6007740
$ws.Range("N68").Value = -6009362
$ws.Range("H71").Value = 2859938
$ws.Range("J71").Value = 2002580
$ws.Range("L71").Value = 18023220
$ws.Range("N71").Value = -18031332
$ws.Range("H108").Value = 2964.5
$ws.Range("I108").Value = 1557.4
$ws.Range("J108").Value = 10000
$ws.Range("K108").Value = 4672.200000000001
$ws.Range("L108").Value = 30000
$ws.Range("M108").Value = -1792.200000000001
$ws.Range("N108").Value = -35760
$ws.Range("H124").Value = 999999
$ws.Range("I124").Value = 999999
$ws.Range("K124").Value = 2999997
$ws.Range("M124").Value = -2995087
$ws.Range("H126").Value = 500499.5
$ws.Range("I126").Value = 500499.5
$ws.Range("K126").Value = 1501498.5
$ws.Range("M126").Value = -1496558.5
$ws.Range("H129").Value = 93655.27
$ws.Range("I129").Value = 955
$ws.Range("J129").Value = 114255.336
$ws.Range("K129").Value = 2865
$ws.Range("L129").Value = 342766.008
$ws.Range("M129").Value = 2135
$ws.Range("N129").Value = -352766.008
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H131").Value = 57500.434
$ws.Range("I131").Value = 201153.6
$ws.Range("J131").Value = 35054.625
$ws.Range("K131").Value = 603460.8
$ws.Range("L131").Value = 105163.875
$ws.Range("M131").Value = -598420.8
$ws.Range("N131").Value = -115243.875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2387.4707
$ws.Range("I102").Value = 1409.3334
$ws.Range("J102").Value = 4735
$ws.Range("K102").Value = 1409.3334
$ws.Range("L102").Value = 4735
$ws.Range("M102").Value = 212.6666
$ws.Range("N102").Value = -7979
$ws.Range("H136").Value = 65645.664
$ws.Range("J136").Value = 65645.664
$ws.Range("L136").Value = 196936.992
$ws.Range("N136").Value = -202036.992
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4470.2256
$ws.Range("J132").Value = 5332.75
$ws.Range("L132").Value = 15998.25
$ws.Range("N132").Value = -21058.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 9999.5
$ws.Range("I61").Value = 9999.5
$ws.Range("K61").Value = 9999.5
$ws.Range("M61").Value = -9707.5
